$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 through 8 (Support Vector Classifier, CART, Random Forest, XGBoost,
# and the old row that used to hold Lasso is handled by updating row 3 below) so
# that only the header row plus the two remaining model rows stay.
$ws.Rows("4:8").Delete()

# Row 2 stays "Logistic Regression" but its metric values changed slightly.
$ws.Range("A2").Value = "Logistic Regression"
$ws.Range("B2").Value = 0.877961432506887
$ws.Range("C2").Value = 0.8787344197087246
$ws.Range("D2").Value = 0.877961432506887
$ws.Range("E2").Value = 0.8766339395517994

# Row 3 used to be "Lasso"; it now becomes "LightGBM" with the metrics that used
# to belong to the LightGBM row.
$ws.Range("A3").Value = "LightGBM"
$ws.Range("B3").Value = 0.8787878787878789
$ws.Range("C3").Value = 0.8801018091608025
$ws.Range("D3").Value = 0.8787878787878789
$ws.Range("E3").Value = 0.8773858443154371
